$wb = $excel.ActiveWorkbook

# --- "Introduction " sheet (history / status table) ---
$wsIntro = $wb.Worksheets.Item("Introduction ")

# New history row: "Update last open point" on 14/2/2020
$wsIntro.Range("E15").Value = "14/2/2020"
$wsIntro.Range("G15").Value = "Update last open point "
$wsIntro.Range("B15").Value = 0.3
$wsIntro.Range("C15").Value = "T.Sharaby"

# "Last update" field gets retyped (typo: missing trailing 0)
$wsIntro.Range("D9").Value = "14/2/202"

# Ref Version bump
$wsIntro.Range("D7").Value = 1.5

# --- "Cross review points " sheet ---
$wsCross = $wb.Worksheets.Item("Cross review points ")

# Last open point resolved
$wsCross.Range("H5").Value = "Resolved"

# --- View / selection state ---
$wsCross.Range("F9").Select()
$wsCross.Application.ActiveWindow.ScrollRow = 1
$wsCross.Application.ActiveWindow.ScrollColumn = 5

$wsIntro.Activate()
$wsIntro.Range("J9").Select()
